$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 31 de Marzo de 2020 a las 06:50"

# Estados Unidos (row 4) updated totals
$ws.Range("D4").Value = 5507
$ws.Range("E4").Value = 155579
$ws.Range("G4").Value = 11
$ws.Range("H4").Value = 3167

# Ecuador (row 33) updated totals
$ws.Range("D33").Value = 54
$ws.Range("E33").Value = 1850

# Birmania (row 157) updated totals
$ws.Range("E157").Value = 13
$ws.Range("G157").Value = 1
$ws.Range("H157").Value = 1

# Bolivia moves up in the country list: remove its current row (114, right
# after Georgia) and re-insert a fresh row right after "Estado de Palestina"
# (row 110), ahead of Camboya, with refreshed case data.
$ws.Rows(114).Delete()
$ws.Rows(111).Insert()
$ws.Range("A111").Value = "Bolivia"
$ws.Range("B111").Value = 107
$ws.Range("C111").Value = 10
$ws.Range("D111").Value = 0
$ws.Range("E111").Value = 101
$ws.Range("F111").Value = 3
$ws.Range("G111").Value = 2
$ws.Range("H111").Value = 6
